# contacts.xlsx: append 12 new submitted-form rows (10-21) to the Contacts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range('A10').Value = "AEDAKULA"
$ws.Range('B10').Value = "'"
$ws.Range('C10').Value = "'"
$ws.Range('D10').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E10').Value = "'8125577708"
$ws.Range('F10').Value = "IT Consulting"
$ws.Range('G10').Value = "'"

# Row 11
$ws.Range('A11').Value = "a.bobby"
$ws.Range('B11').Value = "'"
$ws.Range('C11').Value = "beb"
$ws.Range('D11').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E11').Value = "'8125577708"
$ws.Range('F11').Value = "Managed Services"
$ws.Range('G11').Value = "gggg"

# Row 12
$ws.Range('A12').Value = "bhanu"
$ws.Range('B12').Value = "annna"
$ws.Range('C12').Value = "beb"
$ws.Range('D12').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E12').Value = "'"
$ws.Range('F12').Value = "Mobile Development"
$ws.Range('G12').Value = "jjjjj"

# Row 13
$ws.Range('A13').Value = "Mouli"
$ws.Range('B13').Value = "'"
$ws.Range('C13').Value = "mlrit"
$ws.Range('D13').Value = "rithvikanirvesh0416@gmail.com"
$ws.Range('E13').Value = "'9876543219"
$ws.Range('F13').Value = "Mobile Development"
$ws.Range('G13').Value = "haa"

# Row 14
$ws.Range('A14').Value = "priya"
$ws.Range('B14').Value = "chilukiri"
$ws.Range('C14').Value = "NGO"
$ws.Range('D14').Value = "priya.chilukuri1122@gmail.com"
$ws.Range('E14').Value = "'9550668312"
$ws.Range('F14').Value = "Cloud Services"
$ws.Range('G14').Value = "want an ngo website build for me"

# Row 15
$ws.Range('A15').Value = "AEDAKULA"
$ws.Range('B15').Value = "'"
$ws.Range('C15').Value = "'"
$ws.Range('D15').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E15').Value = "'08125577708"
$ws.Range('F15').Value = "'"
$ws.Range('G15').Value = "'"

# Row 16
$ws.Range('A16').Value = "AEDAKULA"
$ws.Range('B16').Value = "'"
$ws.Range('C16').Value = "'"
$ws.Range('D16').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E16').Value = "'08125577708"
$ws.Range('F16').Value = "'"
$ws.Range('G16').Value = "'"

# Row 17
$ws.Range('A17').Value = "AEDAKULA"
$ws.Range('B17').Value = "'"
$ws.Range('C17').Value = "'"
$ws.Range('D17').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E17').Value = "'08125577708"
$ws.Range('F17').Value = "'"
$ws.Range('G17').Value = "'"

# Row 18
$ws.Range('A18').Value = "AEDAKULA"
$ws.Range('B18').Value = "'"
$ws.Range('C18').Value = "'"
$ws.Range('D18').Value = "rithvikanirvesh2909@gmail.com"
$ws.Range('E18').Value = "'08125577708"
$ws.Range('F18').Value = "'"
$ws.Range('G18').Value = "'"

# Row 19
$ws.Range('A19').Value = "'"
$ws.Range('B19').Value = "'"
$ws.Range('C19').Value = "'"
$ws.Range('D19').Value = "'"
$ws.Range('E19').Value = "'"
$ws.Range('F19').Value = "'"
$ws.Range('G19').Value = "'"

# Row 20
$ws.Range('A20').Value = "'"
$ws.Range('B20').Value = "'"
$ws.Range('C20').Value = "'"
$ws.Range('D20').Value = "'"
$ws.Range('E20').Value = "'"
$ws.Range('F20').Value = "'"
$ws.Range('G20').Value = "'"

# Row 21
$ws.Range('A21').Value = "'"
$ws.Range('B21').Value = "'"
$ws.Range('C21').Value = "'"
$ws.Range('D21').Value = "'"
$ws.Range('E21').Value = "'"
$ws.Range('F21').Value = "'"
$ws.Range('G21').Value = "'"
